# Natmi following Dr Hou advice
# Update the LR-pair (Nts -> Sort1) result table: add sender-cluster "ECs" rows
# alongside the existing "sCs" rows (rows 2-5 = ECs sender, rows 6-9 = sCs sender),
# each cycling through target clusters ECs/FAPs/M2/Nts, with refreshed metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nts"
$ws.Range("C2").Value = "Sort1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.593549
$ws.Range("H2").Value = 4.780647
$ws.Range("I2").Value = 0.8185302145731395
$ws.Range("J2").Value = 0.8185302145731395
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024839333333334
$ws.Range("N2").Value = 9.074518000000001
$ws.Range("O2").Value = 0.1801507982970389
$ws.Range("P2").Value = 0.1801507982970388
$ws.Range("Q2").Value = 4.820229694794001
$ws.Range("R2").Value = 43.382067253146
$ws.Range("S2").Value = 0.1474588715855976
$ws.Range("T2").Value = 0.1474588715855976

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nts"
$ws.Range("C3").Value = "Sort1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.593549
$ws.Range("H3").Value = 4.780647
$ws.Range("I3").Value = 0.8185302145731395
$ws.Range("J3").Value = 0.8185302145731395
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.09684827751501936
$ws.Range("P3").Value = 0.09684827751501934
$ws.Range("Q3").Value = 2.591334302043
$ws.Range("R3").Value = 23.322008718387
$ws.Range("S3").Value = 0.07927324137540775
$ws.Range("T3").Value = 0.07927324137540774

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nts"
$ws.Range("C4").Value = "Sort1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.593549
$ws.Range("H4").Value = 4.780647
$ws.Range("I4").Value = 0.8185302145731395
$ws.Range("J4").Value = 0.8185302145731395
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.640628666666666
$ws.Range("N4").Value = 13.921886
$ws.Range("O4").Value = 0.2763825997921178
$ws.Range("P4").Value = 0.2763825997921177
$ws.Range("Q4").Value = 7.395069171138
$ws.Range("R4").Value = 66.555622540242
$ws.Range("S4").Value = 0.2262275087121243
$ws.Range("T4").Value = 0.2262275087121242

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Nts"
$ws.Range("C5").Value = "Sort1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.593549
$ws.Range("H5").Value = 4.780647
$ws.Range("I5").Value = 0.8185302145731395
$ws.Range("J5").Value = 0.8185302145731395
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.498988000000001
$ws.Range("N5").Value = 22.496964
$ws.Range("O5").Value = 0.4466183243958241
$ws.Range("P5").Value = 0.446618324395824
$ws.Range("Q5").Value = 11.950004828412
$ws.Range("R5").Value = 107.550043455708
$ws.Range("S5").Value = 0.3655705929000099
$ws.Range("T5").Value = 0.3655705929000098

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Nts"
$ws.Range("C6").Value = "Sort1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.353293
$ws.Range("H6").Value = 1.059879
$ws.Range("I6").Value = 0.1814697854268605
$ws.Range("J6").Value = 0.1814697854268605
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.024839333333334
$ws.Range("N6").Value = 9.074518000000001
$ws.Range("O6").Value = 0.1801507982970389
$ws.Range("P6").Value = 0.1801507982970388
$ws.Range("Q6").Value = 1.068654562591334
$ws.Range("R6").Value = 9.617891063322002
$ws.Range("S6").Value = 0.03269192671144127
$ws.Range("T6").Value = 0.03269192671144126

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Nts"
$ws.Range("C7").Value = "Sort1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.353293
$ws.Range("H7").Value = 1.059879
$ws.Range("I7").Value = 0.1814697854268605
$ws.Range("J7").Value = 0.1814697854268605
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.09684827751501936
$ws.Range("P7").Value = 0.09684827751501934
$ws.Range("Q7").Value = 0.5745039967843333
$ws.Range("R7").Value = 5.170535971059
$ws.Range("S7").Value = 0.0175750361396116
$ws.Range("T7").Value = 0.0175750361396116

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nts"
$ws.Range("C8").Value = "Sort1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.353293
$ws.Range("H8").Value = 1.059879
$ws.Range("I8").Value = 0.1814697854268605
$ws.Range("J8").Value = 0.1814697854268605
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.640628666666666
$ws.Range("N8").Value = 13.921886
$ws.Range("O8").Value = 0.2763825997921178
$ws.Range("P8").Value = 0.2763825997921177
$ws.Range("Q8").Value = 1.639501623532667
$ws.Range("R8").Value = 14.755514611794
$ws.Range("S8").Value = 0.05015509107999348
$ws.Range("T8").Value = 0.05015509107999346

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nts"
$ws.Range("C9").Value = "Sort1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.353293
$ws.Range("H9").Value = 1.059879
$ws.Range("I9").Value = 0.1814697854268605
$ws.Range("J9").Value = 0.1814697854268605
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.498988000000001
$ws.Range("N9").Value = 22.496964
$ws.Range("O9").Value = 0.4466183243958241
$ws.Range("P9").Value = 0.446618324395824
$ws.Range("Q9").Value = 2.649339967484
$ws.Range("R9").Value = 23.844059707356
$ws.Range("S9").Value = 0.08104773149581418
$ws.Range("T9").Value = 0.08104773149581417
